$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the FilesTab query (row 4 / cell B4): the "File Type" and
# "Breed" columns are removed from the RETURN clause of the Cypher
# query text stored in that cell.
$newFilesQuery = "`nMATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nWHERE diag.stage_of_disease IN ['IVb']`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFilesQuery

# Update the active sheet view: scroll/selection now sits on the
# FilesTab row (row 4) instead of the CasesTab row (row 2).
$ws.Activate()
$ws.Range("B4").Select()
